# Indicated years and river miles used.
#
# Appends two new paragraphs (matching the existing "Normal (Web)" styled
# paragraphs used for the gage/reach figures) right after the
# "SA gage: 115.96" paragraph and before the document's trailing empty
# paragraph:
#   1. An empty NormalWeb-styled paragraph (spacer).
#   2. A NormalWeb-styled paragraph with the note about the years used.

$d = $word.ActiveDocument

# Locate the "SA gage: 115.96" paragraph so the new content is inserted
# right after it (and therefore before the final, pre-existing empty
# paragraph at the end of the body).
$anchor = $d.Content.Find
$found = $anchor.Execute("SA gage: 115.96", $true, $false, $false, $false,
                          $false, $true, 1, $false, "", 0)

$anchorRange = $d.Content.Find.Parent
$insertionPoint = $anchorRange.End
$rng = $d.Range($insertionPoint, $insertionPoint)

$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$pPr = '<w:pPr><w:pStyle w:val="NormalWeb"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/><w:spacing w:before="0" w:beforeAutospacing="0" w:after="0" w:afterAutospacing="0"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="323130"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr>'
$runPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:cs="Calibri"/><w:color w:val="323130"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr>'

# Empty spacer paragraph.
$blankPara = "<w:p $ns>$pPr</w:p>"

# Paragraph carrying the new note text.
$textPara = "<w:p $ns>$pPr<w:r>$runPr<w:t>Used 2003 to 2018 to make all three datasets consistent</w:t></w:r></w:p>"

# InsertXML merges the trailing "</w:p>" of an inserted fragment with the
# paragraph that currently sits at the insertion point (it doesn't open a
# fresh paragraph there) - so a harmless trailing spacer w:p is appended to
# the fragment to absorb that merge, keeping our two real paragraphs intact
# and the original following paragraph untouched. It is then removed below.
$sentinelPara = "<w:p $ns/>"

$rng.InsertXML($blankPara + $textPara + $sentinelPara)

# Find and remove the now-redundant merged paragraph that resulted from the
# trailing sentinel above (it sits immediately before the document's
# original trailing empty paragraph).
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $p.Range.Start -gt $insertionPoint) {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -eq "`r") {
            $p.Range.Delete()
            break
        }
    }
}
